$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value (quote-prefixed so Excel doesn't
# auto-coerce numeric-/date-looking strings) then strip the resulting
# quote-prefix style back to Normal so no stray style index is left behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$rows = @(
    @{
        row = 48
        A = 112557989; B = 93413; C = "Ovaliderad"; D = "LC"; E = 210
        F = "Grön sköldmossa"; G = "Buxbaumia viridis"
        H = "(Moug. ex Lam. & DC.) Brid. ex Moug. & Nestl."
        I = "3"; J = "kapslar"
        HasL = $true
        P = "Björklund, NV om, Srm"
        Q = 595556; R = 6550792; S = 5
        T = "Södermanland"; U = "Flen"; V = "Södermanland"; W = "Mellösa"
        Y = "2023-08-13"; AA = "2023-08-13"
        AC = "Vid roten på gammal granlåga."
        AD = $false; AE = $false; AG = $false
        AW = "Bo Törnquist"; AX = "Bo Törnquist"
    },
    @{
        row = 49
        A = 112557966; B = 93413; C = "Ovaliderad"; D = "LC"; E = 210
        F = "Grön sköldmossa"; G = "Buxbaumia viridis"
        H = "(Moug. ex Lam. & DC.) Brid. ex Moug. & Nestl."
        I = "3"; J = "kapslar"
        HasL = $true
        P = "Björklund, NV om, Srm"
        Q = 595502; R = 6550929; S = 5
        T = "Södermanland"; U = "Flen"; V = "Södermanland"; W = "Mellösa"
        Y = "2023-08-13"; AA = "2023-08-13"
        AC = "Nära roten på gammal granlåga."
        AD = $false; AE = $false; AG = $false
        AW = "Bo Törnquist"; AX = "Bo Törnquist"
    },
    @{
        row = 50
        A = 112557762; B = 89017; C = "Ovaliderad"; D = "NT"; E = 2008
        F = "Fyrflikig jordstjärna"; G = "Geastrum quadrifidum"
        H = "Pers.:Pers."
        I = "10"; J = "fruktkroppar"
        HasL = $false
        P = "Björklund, NV om, Srm"
        Q = 595561; R = 6550947; S = 5
        T = "Södermanland"; U = "Flen"; V = "Södermanland"; W = "Mellösa"
        Y = "2023-08-13"; AA = "2023-08-13"
        AC = "Grandominerad skog fortfarande oangripen av granbarkborrar."
        AD = $false; AE = $false; AG = $false
        AW = "Bo Törnquist"; AX = "Bo Törnquist"
    }
)

foreach ($r in $rows) {
    $n = $r.row

    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    Set-TextValue $ws.Range("I$n") $r.I
    $ws.Range("J$n").Value = $r.J
    Set-TextValue $ws.Range("K$n") ""
    if ($r.HasL) { Set-TextValue $ws.Range("L$n") "" }
    Set-TextValue $ws.Range("N$n") ""
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
    $ws.Range("U$n").Value = $r.U
    $ws.Range("V$n").Value = $r.V
    $ws.Range("W$n").Value = $r.W
    Set-TextValue $ws.Range("Y$n") $r.Y
    Set-TextValue $ws.Range("AA$n") $r.AA
    $ws.Range("AC$n").Value = $r.AC
    $ws.Range("AD$n").Value = $r.AD
    $ws.Range("AE$n").Value = $r.AE
    Set-TextValue $ws.Range("AF$n") ""
    $ws.Range("AG$n").Value = $r.AG
    Set-TextValue $ws.Range("AT$n") ""
    $ws.Range("AW$n").Value = $r.AW
    $ws.Range("AX$n").Value = $r.AX
    Set-TextValue $ws.Range("AY$n") ""
}
